# This script applies the OOXML diff for
# Code/Results/Cases/Case_5_198/res_bus/vm_pu.xlsx
# It updates columns B,C,D,E,F,I,J,K,L,M for rows 2-25
# (data rows corresponding to index 0..23 in column A)
# with the new per-unit voltage magnitude results for the
# "case with 380 kV" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @('B','C','D','E','F','I','J','K','L','M')
$newData = @(
    ,@(2, @(1.02,1.029585357532214,1.031965529203292,1.02934604942493,1.038377287140346,1.031754744203483,1.034731317919555,1.034772123665849,1.032160227836086,1.041165493558144))
    ,@(3, @(1.02,1.0310575500719,1.033028217522836,1.030613589053387,1.039955187531475,1.032068959609567,1.035841933759486,1.035642805629568,1.033234656961155,1.04255136715764))
    ,@(4, @(1.02,1.0320086449775,1.033714378219901,1.03143268954541,1.040974837827388,1.032270253801178,1.036558705696557,1.036204136795907,1.033928259641451,1.043446286430385))
    ,@(5, @(1.02,1.032408130863589,1.034002492719051,1.03177678581158,1.041403182251292,1.032354395086768,1.036859594230036,1.036439631829918,1.034219466319082,1.043822078907442))
    ,@(6, @(1.02,1.032475185738518,1.034050848163151,1.031834546379265,1.041475084840531,1.032368494519754,1.036910088955412,1.036479143951879,1.034268338837553,1.043885151023348))
    ,@(7, @(1.02,1.03201398431068,1.033718229379802,1.031437288366107,1.040980562620215,1.032271379995986,1.0365627279154,1.036207285406767,1.033932152261474,1.043451309474576))
    ,@(8, @(1.02,1.030083209567404,1.032324976662308,1.029774647846904,1.038910831874038,1.031861354981383,1.035107046407274,1.035066803412009,1.032523675692274,1.041634240150848))
    ,@(9, @(1.02,1.026669015588569,1.029858433108188,1.026836323082768,1.035252939656885,1.031123254302931,1.032527375057896,1.033041176492794,1.030029086613865,1.038417956368884))
    ,@(10, @(1.02,1.02438439114048,1.028206089816871,1.024871346474503,1.032806588049648,1.030620594051657,1.030797447820595,1.031679768205416,1.028357183384409,1.036263635243923))
    ,@(11, @(1.02,1.023393001303178,1.027488650881246,1.024018965605556,1.031745334443221,1.03040039798856,1.03004588567465,1.031087596895346,1.027631060628078,1.035328282748688))
    ,@(12, @(1.02,1.023024425284482,1.02722186160669,1.023702116182853,1.031350831929353,1.030318223319004,1.029766340790519,1.030867231346013,1.027361012803858,1.034980463910071))
    ,@(13, @(1.02,1.02310350116142,1.027279102499815,1.02377009230544,1.031435468053406,1.030335867480889,1.029826321478955,1.030914518995697,1.02741895417286,1.035055089903405))
    ,@(14, @(1.02,1.023362541474316,1.027466604159089,1.02399277961266,1.031712731036809,1.030393613252647,1.030022786224098,1.031069389738522,1.027608745230884,1.03529953988813))
    ,@(15, @(1.02,1.023522100792903,1.02758209013281,1.024129952971578,1.031883521103053,1.030429141368905,1.030143783915055,1.031164756657562,1.027725637382841,1.035450102000656))
    ,@(16, @(1.02,1.02445014069788,1.028253662038787,1.024927883299134,1.032876977647409,1.030635154009346,1.030847273436553,1.031719011911539,1.028405327340489,1.036325657699175))
    ,@(17, @(1.02,1.025031698153301,1.028674391807864,1.025427988626992,1.033499612891702,1.030763698410493,1.031287881795539,1.032065962390352,1.028831091423682,1.036874190328301))
    ,@(18, @(1.02,1.025370705668898,1.028919607302784,1.02571954415783,1.033862596375862,1.030838431186672,1.031544641195684,1.032268075230539,1.029079222732323,1.037193897872349))
    ,@(19, @(1.02,1.025486263961714,1.029003187547873,1.025818932257437,1.03398633259034,1.03086387164562,1.031632149012393,1.032336946966357,1.029163793721981,1.037302869065346))
    ,@(20, @(1.02,1.024969323825083,1.028629271080998,1.025374347359255,1.033432829657382,1.030749932170865,1.031240633562286,1.032028764560327,1.028785432700836,1.036815363092185))
    ,@(21, @(1.02,1.023286269734947,1.027411397934027,1.023927210355002,1.03163109251852,1.03037661918241,1.029964942832469,1.031023795435649,1.027552865752172,1.035227566224446))
    ,@(22, @(1.02,1.022226154661403,1.02664393131335,1.023015961968547,1.030496495272164,1.030139679348135,1.029160656287496,1.030389575451596,1.026775969272902,1.034227010102583))
    ,@(23, @(1.02,1.02278832610357,1.027050946818111,1.023499164717172,1.031098138255298,1.030265497142941,1.029587235492507,1.030726012548531,1.027188002124972,1.034757639707309))
    ,@(24, @(1.02,1.024997508733108,1.028649659779711,1.025398585995033,1.033463006710884,1.030756153299114,1.031261983746907,1.032045573454417,1.028806064560431,1.036841945336744))
    ,@(25, @(1.02,1.027553128538659,1.030497479767423,1.027596998509574,1.036199921239623,1.031315929619931,1.03319604466807,1.033566767144222,1.030675532634671,1.039251194214635))
)

foreach ($entry in $newData) {
    $row = $entry[0]
    $values = $entry[1]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $val = $values[$i]
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $val
    }
}
